# The workbook contains a Markov state-transition probability matrix
# (sheet "Sheet1"): row A holds the "from" state, column 1 the "to" state,
# and each data row sums to 1. More simulated games were added, so the
# transition probabilities for several "from" states need to be refreshed
# with their recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - state "Af0"
$ws.Range("B2").Value = 0.1951951951951952
$ws.Range("C2").Value = 0.5765765765765766
$ws.Range("J2").Value = 0.01201201201201201
$ws.Range("P2").Value = 0.1471471471471471
$ws.Range("S2").Value = 0.06906906906906907

# Row 3 - state "Af1"
$ws.Range("B3").Value = 0.005025125628140704
$ws.Range("C3").Value = 0.03015075376884422
$ws.Range("J3").Value = 0.02010050251256281
$ws.Range("P3").Value = 0.7587939698492462
$ws.Range("S3").Value = 0.185929648241206

# Row 4 - state "Af2"
$ws.Range("J4").Value = 0.02173913043478261
$ws.Range("P4").Value = 0.7173913043478261
$ws.Range("S4").Value = 0.2608695652173913

# Row 5 - state "Af3"
$ws.Range("P5").Value = 0.6666666666666666
$ws.Range("S5").Value = 0.3333333333333333

# Row 6 - state "Ai0"
$ws.Range("B6").Value = 0.08333333333333333
$ws.Range("D6").Value = 0.008771929824561403
$ws.Range("E6").Value = 0.004385964912280702
$ws.Range("F6").Value = 0.09649122807017543
$ws.Range("J6").Value = 0.2192982456140351
$ws.Range("O6").Value = 0.01754385964912281
$ws.Range("Q6").Value = 0.1052631578947368
$ws.Range("R6").Value = 0.04824561403508772
$ws.Range("S6").Value = 0.4166666666666667

# Row 7 - state "Ai1"
$ws.Range("B7").Value = 0.10546875
$ws.Range("D7").Value = 0.015625
$ws.Range("F7").Value = 0.05859375
$ws.Range("J7").Value = 0.15234375
$ws.Range("O7").Value = 0.0234375
$ws.Range("Q7").Value = 0.16015625
$ws.Range("R7").Value = 0.05859375
$ws.Range("S7").Value = 0.42578125

# Row 8 - state "Ai2"
$ws.Range("B8").Value = 0.09715639810426541
$ws.Range("D8").Value = 0.02606635071090047
$ws.Range("E8").Value = 0.002369668246445498
$ws.Range("F8").Value = 0.05213270142180094
$ws.Range("J8").Value = 0.1255924170616114
$ws.Range("O8").Value = 0.03080568720379147
$ws.Range("Q8").Value = 0.1540284360189574
$ws.Range("R8").Value = 0.0924170616113744
$ws.Range("S8").Value = 0.4194312796208531

# Row 9 - state "Ai3"
$ws.Range("B9").Value = 0.09836065573770492
$ws.Range("D9").Value = 0.01639344262295082
$ws.Range("F9").Value = 0.06557377049180328
$ws.Range("J9").Value = 0.08196721311475409
$ws.Range("O9").Value = 0.0273224043715847
$ws.Range("Q9").Value = 0.1693989071038251
$ws.Range("R9").Value = 0.06557377049180328
$ws.Range("S9").Value = 0.4754098360655737

# Row 10 - state "Ar0"
$ws.Range("B10").Value = 0.1223241590214067
$ws.Range("D10").Value = 0.02064220183486239
$ws.Range("E10").Value = 0.001529051987767584
$ws.Range("F10").Value = 0.07339449541284404
$ws.Range("J10").Value = 0.1162079510703364
$ws.Range("O10").Value = 0.03134556574923547
$ws.Range("Q10").Value = 0.1880733944954129
$ws.Range("R10").Value = 0.07874617737003058
$ws.Range("S10").Value = 0.367737003058104

# Row 11 - state "Bf0"
$ws.Range("G11").Value = 0.1381909547738693
$ws.Range("J11").Value = 0.0678391959798995
$ws.Range("K11").Value = 0.1984924623115578
$ws.Range("L11").Value = 0.585427135678392
$ws.Range("S11").Value = 0.01005025125628141

# Row 12 - state "Bf1"
$ws.Range("G12").Value = 0.7366255144032922
$ws.Range("J12").Value = 0.205761316872428
$ws.Range("L12").Value = 0.01234567901234568
$ws.Range("S12").Value = 0.04526748971193416

# Row 13 - state "Bf2"
$ws.Range("G13").Value = 0.6521739130434783
$ws.Range("J13").Value = 0.3260869565217391
$ws.Range("S13").Value = 0.02173913043478261

# Row 15 - state "Bi0"
$ws.Range("F15").Value = 0.02049180327868852
$ws.Range("H15").Value = 0.1147540983606557
$ws.Range("I15").Value = 0.06147540983606557
$ws.Range("J15").Value = 0.3360655737704918
$ws.Range("K15").Value = 0.06557377049180328
$ws.Range("M15").Value = 0.01229508196721311
$ws.Range("O15").Value = 0.06557377049180328
$ws.Range("S15").Value = 0.3237704918032787

# Row 16 - state "Bi1"
$ws.Range("F16").Value = 0.02202643171806168
$ws.Range("H16").Value = 0.1233480176211454
$ws.Range("I16").Value = 0.0881057268722467
$ws.Range("J16").Value = 0.4581497797356828
$ws.Range("K16").Value = 0.1585903083700441
$ws.Range("M16").Value = 0.00881057268722467
$ws.Range("O16").Value = 0.00881057268722467
$ws.Range("S16").Value = 0.13215859030837

# Row 17 - state "Bi2"
$ws.Range("F17").Value = 0.01745635910224439
$ws.Range("H17").Value = 0.1920199501246883
$ws.Range("I17").Value = 0.0598503740648379
$ws.Range("J17").Value = 0.3990024937655861
$ws.Range("K17").Value = 0.1221945137157107
$ws.Range("M17").Value = 0.02992518703241895
$ws.Range("N17").Value = 0.007481296758104738
$ws.Range("O17").Value = 0.04239401496259352
$ws.Range("S17").Value = 0.1296758104738155

# Row 18 - state "Bi3"
$ws.Range("F18").Value = 0.02259887005649718
$ws.Range("H18").Value = 0.192090395480226
$ws.Range("I18").Value = 0.1129943502824859
$ws.Range("J18").Value = 0.3446327683615819
$ws.Range("K18").Value = 0.1525423728813559
$ws.Range("M18").Value = 0.01129943502824859
$ws.Range("O18").Value = 0.05649717514124294
$ws.Range("S18").Value = 0.1073446327683616

# Row 19 - state "Br0"
$ws.Range("F19").Value = 0.01412639405204461
$ws.Range("H19").Value = 0.187360594795539
$ws.Range("I19").Value = 0.07806691449814127
$ws.Range("J19").Value = 0.3776951672862454
$ws.Range("K19").Value = 0.1434944237918216
$ws.Range("M19").Value = 0.02156133828996282
$ws.Range("O19").Value = 0.07434944237918216
$ws.Range("S19").Value = 0.1033457249070632
